$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last data row (87) into two new rows (88 and 89),
# copying formatting/styles and values, then bump the date serial
# in column A for each new row.
$lastRow = 87
$newRows = 2

for ($i = 1; $i -le $newRows; $i++) {
    $dstRow = $lastRow + $i
    $ws.Range("A$lastRow`:J$lastRow").Copy($ws.Range("A$dstRow`:J$dstRow"))
    $ws.Cells.Item($dstRow, 1).Value = $ws.Cells.Item($lastRow, 1).Value2 + $i
}
